# BOT; UPDATE DATA
# Adds the 2020-05-20 (serial 43971) daily row to the "all", "kobe" and
# "other" sheets (pushing the footnote row down by one), and revises the
# last ~16 days of "currently hospitalised" figures on "all"/"kobe" to
# reflect newly confirmed recoveries.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("all")
$ws2 = $wb.Worksheets.Item("kobe")
$ws3 = $wb.Worksheets.Item("other")

# ---------------------------------------------------------------------
# 1) "all": revise columns D (軽症・中等症) / E (中等・軽症) for the last
#    days as patients moved between categories / recovered.
# ---------------------------------------------------------------------
$updates1 = @(
    @("D26",96), @("E26",85), @("D27",96), @("E27",86), @("D28",86), @("E28",76),
    @("D29",84), @("E29",74), @("D30",79), @("E30",69), @("D31",78), @("E31",68),
    @("D32",71), @("E32",61), @("D33",70), @("E33",59), @("D34",66), @("E34",55),
    @("D35",55), @("E35",44), @("D36",55), @("E36",45), @("D37",49), @("E37",39),
    @("D38",46), @("E38",36), @("D39",42), @("E39",35), @("D40",42), @("E40",35),
    @("D41",41), @("E41",34)
)
foreach ($u in $updates1) {
    $ws1.Range($u[0]).Value = $u[1]
}

# ---------------------------------------------------------------------
# 2) "kobe": same revision, columns F / G.
# ---------------------------------------------------------------------
$updates2 = @(
    @("F81",91), @("G81",81), @("F82",91), @("G82",82), @("F83",81), @("G83",72),
    @("F84",79), @("G84",70), @("F85",74), @("G85",65), @("F86",73), @("G86",64),
    @("F87",66), @("G87",57), @("F88",65), @("G88",55), @("F89",61), @("G89",51),
    @("F90",50), @("G90",40), @("F91",50), @("G91",41), @("F92",44), @("G92",35),
    @("F93",41), @("G93",32), @("F94",37), @("G94",31), @("F95",37), @("G95",31),
    @("F96",36), @("G96",30)
)
foreach ($u in $updates2) {
    $ws2.Range($u[0]).Value = $u[1]
}

# ---------------------------------------------------------------------
# 3) Insert the new day's row on each sheet just above the trailing
#    footnote row (Insert() pushes the footnote down and inherits the
#    number formats of the row above, same as the row-above "fill down"
#    Excel performs on a manual row insert).
# ---------------------------------------------------------------------

# --- "all" : new row 43 ---
$ws1.Rows.Item(43).Insert()
$ws1.Range("A43").Value = 43971
$ws1.Range("B43").Value = 283
$ws1.Range("C43").Value = 281
$ws1.Range("D43").Value = 39
$ws1.Range("E43").Value = 34
$ws1.Range("F43").Value = 5
$ws1.Range("G43").Value = 11
$ws1.Range("H43").Value = 231

# --- "kobe" : new row 98 ---
$ws2.Rows.Item(98).Insert()
$ws2.Range("A98").Value = 43971
$ws2.Range("B98").Value = 0
$ws2.Range("C98").Value = 2922
$ws2.Range("D98").Value = 0
$ws2.Range("E98").Value = 283
$ws2.Range("F98").Value = 34
$ws2.Range("G98").Value = 30
$ws2.Range("H98").Value = 4
$ws2.Range("I98").Value = 11
$ws2.Range("J98").Value = 222

# --- "other" : new row 73 ---
$ws3.Rows.Item(73).Insert()
$ws3.Range("A73").Value = 43971
$ws3.Range("B73").Value = 0
$ws3.Range("C73").Value = 14
$ws3.Range("D73").Value = 5
$ws3.Range("E73").Value = 4
$ws3.Range("F73").Value = 1
$ws3.Range("G73").Value = 0
$ws3.Range("H73").Value = 9

# ---------------------------------------------------------------------
# 4) Restore per-sheet selections and make "other" the active tab, as
#    left by whoever ran the update (matches the saved view state).
# ---------------------------------------------------------------------
[void]$ws1.Activate()
$ws1.Range("A43").Select() | Out-Null

[void]$ws2.Activate()
$ws2.Range("F79:J98").Select() | Out-Null

[void]$ws3.Activate()
$ws3.Range("A73").Select() | Out-Null
